# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts per game). Replace the previously-written
# values with the actual strikeout totals.
$kValues = @{
    2  = 3
    3  = 2
    4  = 5
    5  = 10
    6  = 7
    7  = 6
    8  = 1
    9  = 3
    10 = 7
    11 = 3
    12 = 3
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
